# Apply updated values to column F (dSF) on Sheet1, as part of a
# "repull data, push all data, mean calculation" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = 4
    6  = -2
    7  = 1
    8  = 0
    10 = 4
    11 = 0
    12 = -1
    13 = -6
    14 = 2
    16 = -1
    20 = -9
    23 = -7
    24 = -4
    26 = -4
    27 = -1
    28 = -9
    29 = 2
    30 = -4
    31 = -2
    36 = -8
    38 = 0
    39 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
